# Weekly update: insert a new latest-week record at the top of the data
# table (row 259), pushing the existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current first data row of this block (row 259).
# This shifts rows 259:269 down to 260:270 and extends the sheet dimension to T270.
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row with this week's price observation.
$ws.Range("A259").Value = 1
$ws.Range("B259").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C259").Value = "Arica y Parinacota"
$ws.Range("D259").Value = 44783
$ws.Range("E259").Value = 15
$ws.Range("F259").Value = "Fruta"
$ws.Range("G259").Value = 100102
$ws.Range("H259").Value = "Cítricos"
$ws.Range("I259").Value = 100102003
$ws.Range("J259").Value = "Limón"
$ws.Range("K259").Value = "Sin especificar"
$ws.Range("L259").Value = "2a amarillo"
$ws.Range("M259").Value = 200
$ws.Range("N259").Value = 9000
$ws.Range("O259").Value = 10000
$ws.Range("P259").Value = 9500
$ws.Range("Q259").Value = "$/caja 20 kilos"
$ws.Range("R259").Value = "Región de Coquimbo"
$ws.Range("S259").Value = 475
$ws.Range("T259").Value = 20
